$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A78").Value = "cc ordenar tambien por numero de fc para las fcs q estan en la misma fehca"
$ws.Range("B78").Value = "no comenzado"

$ws.Range("A79").Value = "ver calculo de total en reporte de comisiones"
$ws.Range("B79").Value = "no comenzado"

$ws.Range("C67").Select()
